$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15..148 down to 16..149
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new record's data
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 44503
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112005
$ws.Range("G15").Value = "Puerro"
$ws.Range("H15").Value = "Azul de Maquehue"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 65
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 7000
$ws.Range("M15").Value = 7000
$ws.Range("N15").Value = "$/docena de paquetes"
$ws.Range("O15").Value = "Provincia de Cautín"
$ws.Range("P15").Value = 583
$ws.Range("Q15").Value = 12
$ws.Range("R15").Value = "Hortaliza"
